$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2025-12-31'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Wednesday'
$ws.Cells.Item($row, 3).Value = 'Harali KH'
$ws.Cells.Item($row, 4).Value = 'sakshi'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 60
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 'Half paid'
$ws.Cells.Item($row, 12).Value = 30

$row = 14
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2025-12-31'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Wednesday'
$ws.Cells.Item($row, 3).Value = 'vairgwadi'
$ws.Cells.Item($row, 4).Value = 'Suresh Patil'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 60
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 'Not paid'
$ws.Cells.Item($row, 12).Value = 0

$row = 15
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-01-14'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Wednesday'
$ws.Cells.Item($row, 3).Value = 'Harali KH'
$ws.Cells.Item($row, 4).Value = 'sakshi'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 60
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 60
$ws.Cells.Item($row, 11).Value = 'Half paid'
$ws.Cells.Item($row, 12).Value = 11

$row = 16
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2025-12-31'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Wednesday'
$ws.Cells.Item($row, 3).Value = 'vairgwadi'
$ws.Cells.Item($row, 4).Value = 'Anil Dhotare'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 35
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 35
$ws.Cells.Item($row, 11).Value = 'Paid'
$ws.Cells.Item($row, 12).Value = 0

$row = 17
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2025-12-31'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Wednesday'
$ws.Cells.Item($row, 3).Value = 'vairgwadi'
$ws.Cells.Item($row, 4).Value = 'geeta Morti'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 35
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 35
$ws.Cells.Item($row, 11).Value = 'Paid'
$ws.Cells.Item($row, 12).Value = 0

$row = 18
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2025-12-31'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Monday'
$ws.Cells.Item($row, 3).Value = 'Harali KH'
$ws.Cells.Item($row, 4).Value = 'sakshi'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '500gm'
$ws.Cells.Item($row, 8).Value = 170
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 170
$ws.Cells.Item($row, 11).Value = 'Half paid'
$ws.Cells.Item($row, 12).Value = 50

$row = 19
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-01-02'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Friday'
$ws.Cells.Item($row, 3).Value = 'Harali BK'
$ws.Cells.Item($row, 4).Value = 'Ramdas Salve'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '500gm'
$ws.Cells.Item($row, 8).Value = 170
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 170
$ws.Cells.Item($row, 11).Value = 'Half paid'
$ws.Cells.Item($row, 12).Value = 100

$row = 20
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-01-01'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Thursday'
$ws.Cells.Item($row, 3).Value = 'vairgwadi'
$ws.Cells.Item($row, 4).Value = 'Hausabai Murkute'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 35
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 70
$ws.Cells.Item($row, 11).Value = 'Paid'
$ws.Cells.Item($row, 12).Value = 0

$row = 21
$dateCell = $ws.Cells.Item($row, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = '2026-01-01'
$dateCell.Style = "Normal"
$ws.Cells.Item($row, 2).Value = 'Thursday'
$ws.Cells.Item($row, 3).Value = 'Harali BK'
$ws.Cells.Item($row, 4).Value = 'Hari Patake'
$ws.Cells.Item($row, 5).Value = 'GOLD Tea Powder'
$ws.Cells.Item($row, 6).Value = 'Mix'
$ws.Cells.Item($row, 7).Value = '100gm'
$ws.Cells.Item($row, 8).Value = 35
$ws.Cells.Item($row, 9).Value = 1
$ws.Cells.Item($row, 10).Value = 35
$ws.Cells.Item($row, 11).Value = 'Paid'
$ws.Cells.Item($row, 12).Value = 0

